$wb = $excel.ActiveWorkbook
$wsErrores = $wb.Worksheets.Item("Errores")

# Insert two new rows before row 15 of "Errores", shifting the existing
# AuctionManagementBean / BussinessException rows from 15/17 down to 17/19.
[void]$wsErrores.Rows.Item(14).Resize(2, 1).EntireRow.Insert()

# Populate the two new rows with the new comment strings (column C only).
$wsErrores.Range("C14").Value = "Los campos de fecha en la base de datos estan como Date, lo cual impide que se almacene la hora, se cambian a datetime"
$wsErrores.Range("C15").Value = "El mapeo de hibernate tiene los campos fecha como DATE se cambian a TIMESTAMP"

# Move the selection / active sheet to "Errores" (it becomes the active tab).
[void]$wsErrores.Range("C16").Select()

Write-Host "done"
